# Auto-generated Excel COM-interop edit script
# Applies data value corrections + column width tweaks per the commit
# "running pipeline on new warning sets"

$wb = $excel.ActiveWorkbook

# --- Column width adjustments ---
$wsAllTools = $wb.Worksheets.Item("all_tools")
$wsAllTools.Columns.Item(9).ColumnWidth = 20.75   # -> 21.7109375 (raw OOXML width)
$wsAllTools.Columns.Item(11).ColumnWidth = 20.75  # -> 21.7109375 (raw OOXML width)
$wsAllTools.Columns.Item(12).ColumnWidth = 19.75  # -> 20.7109375 (raw OOXML width)

$wsInfer = $wb.Worksheets.Item("infer")
$wsInfer.Columns.Item(12).ColumnWidth = 18.75     # -> 19.7109375 (raw OOXML width)

# --- Sheet: all_tools ---
$ws = $wb.Worksheets.Item("all_tools")

# Value corrections
$ws.Range("L5").Value = 0.3140610795011265
$ws.Range("F9").Value = 91
$ws.Range("G9").Value = 427
$ws.Range("I9").Value = -0.1910649993965447
$ws.Range("J9").Value = 0.007362964502577038
$ws.Range("K9").Value = -0.2759146301993454
$ws.Range("L9").Value = 0.005459194924523375
$ws.Range("F10").Value = 46
$ws.Range("G10").Value = 790
$ws.Range("I10").Value = -0.06118096042217638
$ws.Range("J10").Value = 0.5577933170131231
$ws.Range("K10").Value = -0.0733097603574439
$ws.Range("L10").Value = 0.6128926682147983
$ws.Range("F11").Value = 46
$ws.Range("G11").Value = 790
$ws.Range("I11").Value = -0.04186119177310581
$ws.Range("J11").Value = 0.6748774434313816
$ws.Range("K11").Value = -0.05015044700647268
$ws.Range("L11").Value = 0.7294436579599445
$ws.Range("F12").Value = 46
$ws.Range("G12").Value = 790
$ws.Range("I12").Value = 0.1172265866272641
$ws.Range("J12").Value = 0.237157137469037
$ws.Range("K12").Value = 0.1738547627439666
$ws.Range("L12").Value = 0.2272525783888778
$ws.Range("G13").Value = 65
$ws.Range("I13").Value = -0.3373495424699933
$ws.Range("J13").Value = 0.2074202127647988
$ws.Range("K13").Value = -0.4296689244236597
$ws.Range("L13").Value = 0.215243543278886
$ws.Range("G14").Value = 65
$ws.Range("I14").Value = -0.1816497536376887
$ws.Range("J14").Value = 0.4972433060612282
$ws.Range("K14").Value = -0.2669155439601523
$ws.Range("L14").Value = 0.4559719917038285
$ws.Range("G15").Value = 65
$ws.Range("I15").Value = -0.0524863881081478
$ws.Range("J15").Value = 0.8456867367859529
$ws.Range("K15").Value = -0.06856450678985078
$ws.Range("L15").Value = 0.8507182473580949
$ws.Range("G16").Value = 65
$ws.Range("I16").Value = -0.7525489793561388
$ws.Range("J16").Value = 0.004918698145511134
$ws.Range("K16").Value = -0.8788682545029405
$ws.Range("L16").Value = 0.000811787483996615
$ws.Range("G17").Value = 65
$ws.Range("I17").Value = -0.1816497536376887
$ws.Range("J17").Value = 0.4972433060612282
$ws.Range("K17").Value = -0.260405408741612
$ws.Range("L17").Value = 0.4674445466605421
$ws.Range("G18").Value = 65
$ws.Range("I18").Value = 0.4411494016915297
$ws.Range("J18").Value = 0.09923045565594253
$ws.Range("K18").Value = 0.5143006822646836
$ws.Range("L18").Value = 0.1282920587230653
$ws.Range("G19").Value = 65
$ws.Range("I19").Value = 0.4198911048651824
$ws.Range("J19").Value = 0.1194709867717007
$ws.Range("K19").Value = 0.5060713596393749
$ws.Range("L19").Value = 0.1355782583455031
$ws.Range("G20").Value = 65
$ws.Range("I20").Value = -0.3892494720807615
$ws.Range("J20").Value = 0.1457680056362324
$ws.Range("K20").Value = -0.5077905470461433
$ws.Range("L20").Value = 0.1340355823255553
$ws.Range("G21").Value = 65
$ws.Range("I21").Value = -0.3373495424699933
$ws.Range("J21").Value = 0.2074202127647988
$ws.Range("K21").Value = -0.4687297357349016
$ws.Range("L21").Value = 0.1717865787289185
$ws.Range("G22").Value = 65
$ws.Range("I22").Value = 0.0778498944161523
$ws.Range("J22").Value = 0.7711058640185235
$ws.Range("K22").Value = 0.09765202827810447
$ws.Range("L22").Value = 0.788411563708648
$ws.Range("G23").Value = 65
$ws.Range("I23").Value = 0.1297498240269205
$ws.Range("J23").Value = 0.6277606629910362
$ws.Range("K23").Value = 0.2083243269932896
$ws.Range("L23").Value = 0.5635582121900502
$ws.Range("G24").Value = 65
$ws.Range("I24").Value = -0.2335496832484569
$ws.Range("J24").Value = 0.3827797056047885
$ws.Range("K24").Value = -0.3320168961455552
$ws.Range("L24").Value = 0.3486190102393061

# --- Sheet: checker_framework ---
$ws = $wb.Worksheets.Item("checker_framework")

# Value corrections
$ws.Range("F9").Value = 18
$ws.Range("G9").Value = 51
$ws.Range("I9").Value = -0.2419880147043038
$ws.Range("J9").Value = 0.002471409738675847
$ws.Range("K9").Value = -0.303489343721962
$ws.Range("L9").Value = 0.002144343888314521
$ws.Range("L11").Value = 0.00491305256761129

# --- Sheet: typestate_checker ---
$ws = $wb.Worksheets.Item("typestate_checker")

# Value corrections
$ws.Range("L14").Value = 0.4230203924441357

# --- Sheet: infer ---
$ws = $wb.Worksheets.Item("infer")

# Value corrections
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5
$ws.Range("I10").Value = 0.01421997815784012
$ws.Range("J10").Value = 0.9088173646349219
$ws.Range("K10").Value = 0.01636148293791983
$ws.Range("L10").Value = 0.9102093962355748
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 5
$ws.Range("I11").Value = 0.001918117792161453
$ws.Range("J11").Value = 0.9870933347767511
$ws.Range("K11").Value = 0.002310973073108823
$ws.Range("L11").Value = 0.9872920353883015
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 5
$ws.Range("I12").Value = 0.2419047619047619
$ws.Range("J12").Value = 0.04002718189621234
$ws.Range("K12").Value = 0.2933526131391836
$ws.Range("L12").Value = 0.03867934687031337
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0

# Clear now-undefined correlation stats (0/0 sample sizes)
$ws.Range("I9:L9").ClearContents()
$ws.Range("I13:L13").ClearContents()
$ws.Range("I14:L14").ClearContents()
$ws.Range("I15:L15").ClearContents()
$ws.Range("I16:L16").ClearContents()
$ws.Range("I17:L17").ClearContents()
$ws.Range("I18:L18").ClearContents()
$ws.Range("I19:L19").ClearContents()
$ws.Range("I20:L20").ClearContents()
$ws.Range("I21:L21").ClearContents()
$ws.Range("I22:L22").ClearContents()
$ws.Range("I23:L23").ClearContents()
$ws.Range("I24:L24").ClearContents()

# --- Sheet: openjml ---
$ws = $wb.Worksheets.Item("openjml")

# Value corrections
$ws.Range("L7").Value = 0.0604349562009266

